$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("IMX-USD", "TAO-USD", "GRT-USD", "PEPE-USD", "MNT-USD")

$startRow = 379
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
